$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "cryptos" price/volume snapshot to the newest scrape.
# Note: several "Price" column values look numeric (e.g. "1.000", "0.1520")
# but must stay as literal text (leading zeros / trailing zeros matter,
# and some use "." as a thousands separator). A leading apostrophe forces
# Excel to store them as text instead of re-parsing them as numbers.

$ws.Range("D2").Value = '27.779.77'
$ws.Range("E2").Value = '  -1.60%  '
$ws.Range("D3").Value = '1.894.89'
$ws.Range("E3").Value = '  -1.40%  '
$ws.Range("E4").Value = '  -0.61%  '
$ws.Range("D5").Value = '''311.94'
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("D7").Value = '''0.4945'
$ws.Range("E7").Value = '  +1.88%  '
$ws.Range("D8").Value = '''0.3794'
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").Value = '''0.07324'
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("E10").Value = '  -3.82%  '
$ws.Range("D11").Value = '''20.60'
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").Value = '''0.07623'
$ws.Range("E12").Value = '  -2.22%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''5.460'
$ws.Range("E13").Value = '  -1.58%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.828.45'
$ws.Range("E14").Value = '  -5.11%  '
$ws.Range("D15").Value = '''6.644'
$ws.Range("D16").Value = '''91.05'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '''0.000008733'
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '''0.9998'
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("D20").Value = '27.853.94'
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("E21").Value = '  -3.50%  '
$ws.Range("D22").Value = '''5.116'
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").Value = '2.107.76'
$ws.Range("E23").Value = '  -2.50%  '
$ws.Range("D24").Value = '''10.75'
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("D25").Value = '''154.05'
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("D26").Value = '''1.851'
$ws.Range("E26").Value = '  -4.28%  '
$ws.Range("D27").Value = '''2.186'
$ws.Range("E27").Value = '  +3.71%  '
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("D29").Value = '''115.08'
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").Value = '''4.881'
$ws.Range("E30").Value = '  -2.69%  '
$ws.Range("D31").Value = '''0.08940'
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").Value = '''3.252'
$ws.Range("E32").Value = '  -3.33%  '
$ws.Range("D33").Value = '''1.229'
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("D34").Value = '''0.7671'
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").Value = '''4.644'
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("D36").Value = '''0.02049'
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("D37").Value = '''2.555'
$ws.Range("E37").Value = '  -7.60%  '
$ws.Range("E38").Value = '  -2.82%  '
$ws.Range("D39").Value = '''0.5501'
$ws.Range("E39").Value = '  -1.39%  '
$ws.Range("D40").Value = '''0.05284'
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("D41").Value = '''2.989'
$ws.Range("E41").Value = '  -1.74%  '
$ws.Range("D42").Value = '''6.904'
$ws.Range("E42").Value = '  -2.89%  '
$ws.Range("D43").Value = '''8.543'
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''112.69'
$ws.Range("E44").Value = '  +4.99%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '''0.1520'
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").Value = '''10.56'
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("D47").Value = '''0.4789'
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("D48").Value = '''0.9998'
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("D49").Value = '''1.632'
$ws.Range("E49").Value = '  -2.84%  '
$ws.Range("D50").Value = '''67.45'
$ws.Range("E50").Value = '  -3.09%  '
$ws.Range("D51").Value = '''0.06060'
$ws.Range("E51").Value = '  -1.58%  '
